$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.797.34"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.941.03"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'552.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'132.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.64%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.88%  "
$ws.Range("D9").Value = "2.936.83"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "'4.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "'32.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "3.429.82"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "'6.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.65%  "
$ws.Range("D18").Value = "2.942.61"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "57.806.37"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "'416.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("E22").Value = "  +7.97%  "
$ws.Range("D23").Value = "'13.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.25%  "
$ws.Range("D24").Value = "'7.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "'78.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "'2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.48%  "
$ws.Range("E30").Value = "  +5.59%  "
$ws.Range("D31").Value = "'25.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("D32").Value = "'5.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'0.0968"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("E34").Value = "  +6.63%  "
$ws.Range("E35").Value = "  +6.39%  "
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("B37").Value = "Cosmos"
$ws.Range("C37").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D37").Value = "'8.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.49%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0698"
$ws.Range("E38").Value = "  +14.17%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'48.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'2.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.72%  "
$ws.Range("D41").Value = "'380.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.52%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "2.702.42"
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'124.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "'1.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "'22.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("E51").Value = "  +4.18%  "
